# CampaignData.xlsx maintenance update
#
#  - "New Classes": the automation's trigger-campaign test fixture used to
#    reuse the broadcast campaign's name ("CampaignSMSBrod19may"); it now
#    gets its own dedicated value ("CampaignSMSTrig19may") on the
#    CamapignTrigger sheet.
#  - The workbook was last saved with the Credentials sheet/tab active and
#    selected; it is now saved with CampaignSMSBroadcast as the active
#    sheet/tab instead.
#  - The Credentials sheet picks up a blank, formatted D2 cell (matching
#    the already-defined column D styling) which extends its used range.

$wb = $excel.ActiveWorkbook

# --- New Classes: update the CamapignTrigger sheet's campaign name ---
$wsTrigger = $wb.Worksheets.Item("CamapignTrigger")
$wsTrigger.Range("A1").Value = "CampaignSMSTrig19may"

# --- Credentials sheet: extend used range with a formatted blank D2 cell ---
$wsCred = $wb.Worksheets.Item("Credentials")
$wsCred.Range("D2").Style = $wsCred.Range("B1").Style

# --- Changed Excel ... Location: active tab moves from Credentials to
#     CampaignSMSBroadcast ---
$wsBroadcast = $wb.Worksheets.Item("CampaignSMSBroadcast")
$wsBroadcast.Activate()
